$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so they are not
# auto-converted to numbers by Excel (mirrors the original text values).
$textForceCells = @(
    'D5', 'D6', 'D12', 'D13', 'D15', 'D19', 'D20', 'D22',
    'D23', 'D24', 'D25', 'D27', 'D28', 'D29', 'D32', 'D33',
    'D35', 'D36', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43',
    'D44', 'D45', 'D46', 'D47', 'D48', 'D51'
)
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply updated cell values per the source diff
$ws.Range('D2').Value = '58.646.37'
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').Value = '2.488.10'
$ws.Range('E3').Value = '  -1.24%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '532.07'
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('D6').Value = '135.71'
$ws.Range('E6').Value = '  -3.19%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '2.508.85'
$ws.Range('E9').Value = '  -0.66%  '
$ws.Range('E10').Value = '  +1.26%  '
$ws.Range('E11').Value = '  -0.67%  '
$ws.Range('D12').Value = '5.30'
$ws.Range('E12').Value = '  -1.75%  '
$ws.Range('D13').Value = '0.346'
$ws.Range('E13').Value = '  -2.37%  '
$ws.Range('D14').Value = '2.944.72'
$ws.Range('E14').Value = '  -0.77%  '
$ws.Range('D15').Value = '23.01'
$ws.Range('E15').Value = '  -0.56%  '
$ws.Range('D16').Value = '58.614.51'
$ws.Range('E16').Value = '  -0.87%  '
$ws.Range('E17').Value = '  -1.57%  '
$ws.Range('D18').Value = '2.502.97'
$ws.Range('E18').Value = '  -0.63%  '
$ws.Range('D19').Value = '10.99'
$ws.Range('E19').Value = '  +0.30%  '
$ws.Range('D20').Value = '324.95'
$ws.Range('E20').Value = '  +1.46%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').Value = '5.83'
$ws.Range('E23').Value = '  +0.60%  '
$ws.Range('D24').Value = '63.79'
$ws.Range('E24').Value = '  +2.17%  '
$ws.Range('D25').Value = '0.416'
$ws.Range('E25').Value = '  -0.53%  '
$ws.Range('E26').Value = '  -1.22%  '
$ws.Range('D27').Value = '0.995'
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('D28').Value = '7.54'
$ws.Range('E28').Value = '  -3.38%  '
$ws.Range('D29').Value = '6.70'
$ws.Range('E29').Value = '  -0.77%  '
$ws.Range('D30').Value = '0.0₃0766'
$ws.Range('E30').Value = '  -0.36%  '
$ws.Range('E31').Value = '  -2.05%  '
$ws.Range('D32').Value = '167.06'
$ws.Range('E32').Value = '  +3.28%  '
$ws.Range('D33').Value = '1.15'
$ws.Range('E33').Value = '  +1.61%  '
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '1.37'
$ws.Range('E35').Value = '  -5.16%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').Value = '18.43'
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('E37').Value = '  -3.15%  '
$ws.Range('D38').Value = '1.56'
$ws.Range('E38').Value = '  -1.73%  '
$ws.Range('D39').Value = '36.62'
$ws.Range('E39').Value = '  -0.83%  '
$ws.Range('D40').Value = '0.813'
$ws.Range('E40').Value = '  +1.29%  '
$ws.Range('D41').Value = '3.60'
$ws.Range('E41').Value = '  -0.94%  '
$ws.Range('D42').Value = '5.22'
$ws.Range('E42').Value = '  -1.75%  '
$ws.Range('D43').Value = '277.75'
$ws.Range('E43').Value = '  -2.93%  '
$ws.Range('D44').Value = '0.996'
$ws.Range('E44').Value = '  -0.21%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').Value = '0.598'
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').Value = '10.86'
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').Value = '126.66'
$ws.Range('E47').Value = '  +1.64%  '
$ws.Range('D48').Value = '0.0924'
$ws.Range('E48').Value = '  -0.46%  '
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('E50').Value = '  -0.91%  '
$ws.Range('D51').Value = '17.27'
$ws.Range('E51').Value = '  -1.23%  '

Write-Host "Applied cryptos update."
